$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the existing data rows (2-4) into new rows (5-7), replicating
# the CSV importer test fixture with repeated event-attribute rows.
for ($i = 0; $i -lt 3; $i++) {
    $srcRow = 2 + $i
    $dstRow = 5 + $i
    for ($col = 1; $col -le 5; $col++) {
        $srcCell = $ws.Cells.Item($srcRow, $col)
        $dstCell = $ws.Cells.Item($dstRow, $col)
        $dstCell.Value = $srcCell.Value2
    }
}

# Update the selection to reflect where the cursor ended up after the edit.
$ws.Range("A5").Select()
